$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Debit")
$ws.Range("B2").Value = "Tue Apr 25 14:22:13 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:22:27 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:22:41 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:22:56 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:23:10 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:23:24 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:23:39 EDT 2023"

$ws = $wb.Worksheets.Item("Debit-ZeroDollar")
$ws.Range("B2").Value = "Tue Apr 25 14:23:53 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:24:07 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:24:21 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:24:35 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:24:49 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:25:04 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:25:18 EDT 2023"

$ws = $wb.Worksheets.Item("Debit-Void")
$ws.Range("B2").Value = "Tue Apr 25 14:25:33 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:25:56 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:26:20 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:26:44 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:27:07 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:27:31 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:27:54 EDT 2023"

$ws = $wb.Worksheets.Item("Debit-Credit")
$ws.Range("B2").Value = "Tue Apr 25 14:28:19 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:28:43 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:29:08 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:29:33 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:29:58 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:30:23 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:30:49 EDT 2023"

$ws = $wb.Worksheets.Item("Debit-Credit-Void")
$ws.Range("B2").Value = "Tue Apr 25 14:31:14 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:31:48 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:32:22 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:32:56 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:33:30 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:34:05 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:34:38 EDT 2023"

$ws = $wb.Worksheets.Item("DebitCredit-RemID-Pipe")
$ws.Range("B2").Value = "Tue Apr 25 14:35:13 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:35:38 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:36:03 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:36:28 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:36:53 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:37:17 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:37:43 EDT 2023"

$ws = $wb.Worksheets.Item("Debit-RemID-Pipe")
$ws.Range("B2").Value = "Tue Apr 25 14:38:09 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:38:23 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:38:38 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:38:52 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:39:06 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:39:20 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:39:34 EDT 2023"

$ws = $wb.Worksheets.Item("DebitVoid-RemID-Pipe")
$ws.Range("B2").Value = "Tue Apr 25 14:39:49 EDT 2023"
$ws.Range("B3").Value = "Tue Apr 25 14:40:13 EDT 2023"
$ws.Range("B4").Value = "Tue Apr 25 14:40:37 EDT 2023"
$ws.Range("B5").Value = "Tue Apr 25 14:41:00 EDT 2023"
$ws.Range("B6").Value = "Tue Apr 25 14:41:24 EDT 2023"
$ws.Range("B7").Value = "Tue Apr 25 14:41:48 EDT 2023"
$ws.Range("B8").Value = "Tue Apr 25 14:42:13 EDT 2023"

